$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 (shifts existing rows 10-76 down to 11-77)
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new record
$ws.Cells.Item(10, 1).Value = 2
$ws.Cells.Item(10, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(10, 3).Value = "Coquimbo"
$ws.Cells.Item(10, 4).Value = "2023-05-04"
$ws.Cells.Item(10, 5).Value = 4
$ws.Cells.Item(10, 6).Value = 100112032
$ws.Cells.Item(10, 7).Value = "Zapallo italiano"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 600
$ws.Cells.Item(10, 11).Value = 9500
$ws.Cells.Item(10, 12).Value = 10000
$ws.Cells.Item(10, 13).Value = 9750
$ws.Cells.Item(10, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(10, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(10, 16).Value = 162
$ws.Cells.Item(10, 17).Value = 60
$ws.Cells.Item(10, 18).Value = "Hortaliza"

$ws.Cells.Item(10, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
